$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.616372666666667
$ws.Range("H2").Value = 13.849118
$ws.Range("I2").Value = 0.0183283362562958
$ws.Range("J2").Value = 0.01832833625629581
$ws.Range("M2").Value = 2.843949
$ws.Range("N2").Value = 8.531846999999999
$ws.Range("O2").Value = 0.4976240243095911
$ws.Range("P2").Value = 0.4976240243095912
$ws.Range("Q2").Value = 13.128728428994
$ws.Range("R2").Value = 118.158555860946
$ws.Range("S2").Value = 0.009120620446757302
$ws.Range("T2").Value = 0.009120620446757306
$ws.Range("G3").Value = 4.616372666666667
$ws.Range("H3").Value = 13.849118
$ws.Range("I3").Value = 0.0183283362562958
$ws.Range("J3").Value = 0.01832833625629581
$ws.Range("O3").Value = 0.4403664892852895
$ws.Range("P3").Value = 0.4403664892852897
$ws.Range("Q3").Value = 11.61811280128067
$ws.Range("R3").Value = 104.563015211526
$ws.Range("S3").Value = 0.00807118509162527
$ws.Range("T3").Value = 0.008071185091625273
$ws.Range("G4").Value = 4.616372666666667
$ws.Range("H4").Value = 13.849118
$ws.Range("I4").Value = 0.0183283362562958
$ws.Range("J4").Value = 0.01832833625629581
$ws.Range("M4").Value = 0.3543876666666666
$ws.Range("O4").Value = 0.06200948640511928
$ws.Range("P4").Value = 0.0620094864051193
$ws.Range("Q4").Value = 1.635985537803778
$ws.Range("R4").Value = 14.723869840234
$ws.Range("S4").Value = 0.001136530717913229
$ws.Range("T4").Value = 0.00113653071791323
$ws.Range("I5").Value = 0.943783113604627
$ws.Range("J5").Value = 0.9437831136046271
$ws.Range("M5").Value = 2.843949
$ws.Range("N5").Value = 8.531846999999999
$ws.Range("O5").Value = 0.4976240243095911
$ws.Range("P5").Value = 0.4976240243095912
$ws.Range("Q5").Value = 676.039113486328
$ws.Range("R5").Value = 6084.352021376952
$ws.Range("S5").Value = 0.4696491510673704
$ws.Range("T5").Value = 0.4696491510673706
$ws.Range("I6").Value = 0.943783113604627
$ws.Range("J6").Value = 0.9437831136046271
$ws.Range("O6").Value = 0.4403664892852895
$ws.Range("P6").Value = 0.4403664892852897
$ws.Range("S6").Value = 0.4156104563848091
$ws.Range("T6").Value = 0.4156104563848093
$ws.Range("I7").Value = 0.943783113604627
$ws.Range("J7").Value = 0.9437831136046271
$ws.Range("M7").Value = 0.3543876666666666
$ws.Range("O7").Value = 0.06200948640511928
$ws.Range("P7").Value = 0.0620094864051193
$ws.Range("Q7").Value = 84.24199027613422
$ws.Range("R7").Value = 758.1779124852079
$ws.Range("S7").Value = 0.05852350615244726
$ws.Range("T7").Value = 0.05852350615244729
$ws.Range("I8").Value = 0.03788855013907712
$ws.Range("J8").Value = 0.03788855013907712
$ws.Range("M8").Value = 2.843949
$ws.Range("N8").Value = 8.531846999999999
$ws.Range("O8").Value = 0.4976240243095911
$ws.Range("P8").Value = 0.4976240243095912
$ws.Range("Q8").Value = 27.139860290014
$ws.Range("R8").Value = 244.258742610126
$ws.Range("S8").Value = 0.01885425279546327
$ws.Range("T8").Value = 0.01885425279546328
$ws.Range("I9").Value = 0.03788855013907712
$ws.Range("J9").Value = 0.03788855013907712
$ws.Range("O9").Value = 0.4403664892852895
$ws.Range("P9").Value = 0.4403664892852897
$ws.Range("S9").Value = 0.01668484780885506
$ws.Range("T9").Value = 0.01668484780885507
$ws.Range("I10").Value = 0.03788855013907712
$ws.Range("J10").Value = 0.03788855013907712
$ws.Range("M10").Value = 0.3543876666666666
$ws.Range("O10").Value = 0.06200948640511928
$ws.Range("P10").Value = 0.0620094864051193
$ws.Range("Q10").Value = 3.381928354494888
$ws.Range("S10").Value = 0.002349449534758783
$ws.Range("T10").Value = 0.002349449534758784
